$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# --- Block 1 (Tensor table): insert STDEV row right after the AVG row (row 10) ---
$ws.Rows.Item(11).Insert()
$ws.Range("A10:D10").Copy()
$ws.Range("A11:D11").PasteSpecial($xlPasteFormats)
$ws.Range("A11").Value = "STDEV"
$ws.Range("B11").Formula = "=STDEV(B2:B9)"
$ws.Range("C11").Formula = "=STDEV(C2:C9)"
$ws.Range("D11").Formula = "=STDEV(D2:D9)"

# --- Block 2 (Motion table): insert STDEV row right after its AVG row (originally row 22, now row 23) ---
$ws.Rows.Item(24).Insert()
$ws.Range("A23:D23").Copy()
$ws.Range("A24:D24").PasteSpecial($xlPasteFormats)
$ws.Range("A24").Value = "STDEV"
$ws.Range("B24").Formula = "=STDEV(B15:B22)"
$ws.Range("C24").Formula = "=STDEV(C15:C22)"
$ws.Range("D24").Formula = "=STDEV(D15:D22)"

# --- Block 3 (Images table): insert STDEV row right after its AVG row (originally row 35, now row 37) ---
$ws.Rows.Item(38).Insert()
$ws.Range("A37:D37").Copy()
$ws.Range("A38:D38").PasteSpecial($xlPasteFormats)
$ws.Range("A38").Value = "STDEV"
$ws.Range("B38").Formula = "=STDEV(B29:B36)"
$ws.Range("C38").Formula = "=STDEV(C29:C36)"
$ws.Range("D38").Formula = "=STDEV(D29:D36)"

# --- Restore the selection state to match the authored edit ---
$ws.Range("H23").Select()
